# Update "想去人数" (F column) counts on the 展览, 演出 and 全部类型 sheets
# to reflect the latest generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 83
$ws.Range("F4").Value = 128
$ws.Range("F5").Value = 1750
$ws.Range("F7").Value = 1007
$ws.Range("F8").Value = 2176
$ws.Range("F9").Value = 2090
$ws.Range("F10").Value = 1091
$ws.Range("F11").Value = 595
$ws.Range("F12").Value = 19
$ws.Range("F13").Value = 1656
$ws.Range("F18").Value = 193
$ws.Range("F19").Value = 1558
$ws.Range("F20").Value = 606
$ws.Range("F22").Value = 585
$ws.Range("F23").Value = 12154
$ws.Range("F24").Value = 12183
$ws.Range("F25").Value = 903
$ws.Range("F26").Value = 694
$ws.Range("F28").Value = 23
$ws.Range("F30").Value = 339
$ws.Range("F31").Value = 1910
$ws.Range("F33").Value = 562

# --- 演出 sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 39
$ws.Range("F7").Value = 17

# --- 全部类型 sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 83
$ws.Range("F5").Value = 128
$ws.Range("F6").Value = 1750
$ws.Range("F8").Value = 1007
$ws.Range("F9").Value = 2176
$ws.Range("F10").Value = 2090
$ws.Range("F11").Value = 1091
$ws.Range("F12").Value = 595
$ws.Range("F13").Value = 19
$ws.Range("F14").Value = 1656
$ws.Range("F21").Value = 39
$ws.Range("F22").Value = 193
$ws.Range("F23").Value = 1558
$ws.Range("F24").Value = 606
$ws.Range("F26").Value = 585
$ws.Range("F27").Value = 12154
$ws.Range("F28").Value = 12183
$ws.Range("F29").Value = 903
$ws.Range("F30").Value = 694
$ws.Range("F32").Value = 23
$ws.Range("F34").Value = 339
$ws.Range("F35").Value = 1910
$ws.Range("F39").Value = 562
$ws.Range("F40").Value = 17
